# Auto-generated edit script: update crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.143.45'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '1.863.59'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').Value = '''0.9999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''0.7093'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').Value = '''241.34'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').Value = '''0.07627'
$ws.Range('E9').Value = '  -3.55%  '
$ws.Range('D10').Value = '''24.57'
$ws.Range('E10').Value = '  -3.06%  '
$ws.Range('D11').Value = '''0.08346'
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.865.26'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '''5.209'
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('D14').Value = '''0.7070'
$ws.Range('E14').Value = '  -3.70%  '
$ws.Range('D15').Value = '''91.13'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '29.193.97'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '''242.57'
$ws.Range('E18').Value = '  -2.13%  '
$ws.Range('D19').Value = '''0.000007799'
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('D20').Value = '2.113.15'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('D21').Value = '''13.06'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('D22').Value = '''1.0000'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = '''7.859'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '''0.1586'
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').Value = '''163.92'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').Value = '''8.945'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').Value = '''18.41'
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''1.323'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.499'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('D31').Value = '''4.384'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').Value = '''4.238'
$ws.Range('E32').Value = '  +2.66%  '
$ws.Range('D33').Value = '''0.05144'
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('D34').Value = '''0.7960'
$ws.Range('E34').Value = '  +9.08%  '
$ws.Range('E35').Value = '  -2.49%  '
$ws.Range('D36').Value = '''1.161'
$ws.Range('E36').Value = '  -3.42%  '
$ws.Range('D37').Value = '''2.690'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').Value = '''0.01844'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = '''2.701'
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D40').Value = '1.164.83'
$ws.Range('E40').Value = '  -5.09%  '
$ws.Range('D41').Value = '''6.234'
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').Value = '''0.8897'
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('D43').Value = '''72.80'
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').Value = '2.010.41'
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').Value = '''0.5181'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').Value = '''9.302'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '''0.9996'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = '''0.4267'
$ws.Range('E51').Value = '  -1.79%  '
